# Review_155.docx -> Review 154 "Context is Environment" edit
$d = $word.ActiveDocument

# 1) Heading title
$d.Content.Find.Execute(
    "Review 155: [Short] CHAIN-OF-VERIFICATION REDUCES HALLUCINATION IN LARGE LANGUAGE MODELS, 27.09.2023",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Review 154: Context is Environment, 26.09.2023", 2) | Out-Null

# 2) Bold "Paper:" link
$d.Content.Find.Execute(
    "Paper: https://arxiv.org/abs/2309.11495v2",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Paper: https://arxiv.org/abs/2309.09888v2", 2) | Out-Null

# 3) PDF link line (drop the ".pdf" suffix)
$d.Content.Find.Execute(
    "https://arxiv.org/abs/2309.09888.pdf",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "https://arxiv.org/abs/2309.09888", 2) | Out-Null

# 4) Insert four new paragraphs right after the link paragraph (still index 4).
$linkPara = $d.Paragraphs(4)

$linkPara.Range.InsertParagraphAfter()
$cur = $d.Paragraphs(5)
$cur.Range.Text = "סקירה זו נכתבה על ידי עדן יבין"

$cur.Range.InsertParagraphAfter()
$cur = $d.Paragraphs(6)
# empty paragraph, leave run empty

$cur.Range.InsertParagraphAfter()
$cur = $d.Paragraphs(7)
# empty paragraph, leave run empty

$cur.Range.InsertParagraphAfter()
$cur = $d.Paragraphs(8)
$cur.Range.Text = "מודל של רכב אוטונומי צריך לדעת להתמודד עם המון מצבי עולם אשר לא ראה בזמן תהליך האימון. כיצד הוא עושה זאת? מזעור הטעות על דוגמאות כאלו הינו תחום שלם הנקרא Domain Generalization. האם מודלי שפה יוכלו לעזור לתחום זה ולהראות שיפור על פני המצב הקיים? "

# 5) The paragraph that used to hold only a manual line-break (now paragraph 9)
#    gets three more runs of text appended after the existing <w:br/>. Collapse
#    to the end of the range first so the new text lands in its own <w:t> run
#    (matching Word's own whitespace-preservation rules) instead of forcing
#    xml:space="preserve" the way InsertAfter would.
$brPara = $d.Paragraphs(9)
$tail = $brPara.Range
$tail.Collapse(0)
$tail.Text = "נגלה היום ב- #shorthebrewpapereviews. המאמר נקרא Context is Environments ובמסגרתו החוקרים מנסים להראות שמה שקוראים לו ״הסביבה״ בתחום ה-DG מקביל מאוד ל-Context בתחום של מודלי שפה. שיטות קיימות בתחום ה-DG מנסות להשתמש בדוגמאות העבר שנוצרו מאינטראקציה עם הסביבה כדי לחזות את התוצאה של הדגימה הנוכחית. " + [char]11 + [char]11 + "אך האם לא כך גם אצל מודלי שפה? הרי הם משתמשים בטוקנים הקודמים כדי לחזות את הטוקן הנוכחי. יותר מכך, עם השימוש ההולך וגובר במודלי שפה גילו את היכולת שלהם ללמוד in-context באמצעות טכניקות כגון few-shot. נוכל להשתמש בכך בשביל לשפר את יכולת ההלכה של מודלים אלו על דוגמאות אשר לא ראו. השיטה של החוקרים נקראת ICRM, ובקצרה מנסה להשתמש בקונטקסט כסביבה כדי להקטין את הסיכוי לטעויות על דוגמאות אשר לא נראו ולא דומות למה שהיה באימון המודל. "

# 6) Insert two more paragraphs after that one (now paragraph 9)
$brPara.Range.InsertParagraphAfter()
$cur = $d.Paragraphs(10)
# empty paragraph, leave run empty

$cur.Range.InsertParagraphAfter()
$cur = $d.Paragraphs(11)
$cur.Range.Text = "כאשר מודל השפה h מנסה לשערך את (P (Y| X,C על ידי שימוש בפונקציית הפסד של binary cross-entropy loss. השערוך של (P(Y|X,C) הינו בשביל לשערך את הסיכון של טעות בחיזוי בהינתן הדוגמא הנוכחית והסביבה או הקונטקסט. "

# 7) The "hallucinations" paragraph (now paragraph 12): drop its lead-in text,
#    keep the existing <w:br/>, and append the new closing text after it.
$halluPara = $d.Paragraphs(12)
$find = $halluPara.Range.Duplicate
$find.Find.Execute(
    "מכירים את בעיית ההזיות (hallucinations)במודלי שפה? בגדול זה קורה כאשר מודל שפה מספק לנו תשובות לא נכונות לשאלות לפעמים יחסית פשוטות. סוגיה זו קיבלה התייחסות רבה לאחרונה במספר עבודות ומאמר שנסקור היום ב-#shorthebrewpapereviews מציע גישה נוספת לפתרונה. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "", 2) | Out-Null

$halluTail = $halluPara.Range
$halluTail.Collapse(0)
$halluTail.Text = "החוקרים מראים ששימוש פשוט זה מביא לתוצאות טובות יותר מהשיטות הקודמות בניסויים הכוללים יכולת הכללה על דוגמאות חדשות שלא נראו בסט האימון. למי שירצה להתעמק יותר, המאמר מראה עוד המון נקודות קריטיות וחשובות בשימוש של מודלי שפה בשביל לחשב סיכון של דוגמאות חדשות ובנוסף נותן עוד תאוריה מעניינת על התחום."

# 8) Remove the six now-obsolete paragraphs that followed (the "idea behind
#    the method" paragraph plus the 5-step walkthrough list), which the new
#    review text doesn't need.
$startPara = $d.Paragraphs(13)
$endPara = $d.Paragraphs(18)
$deadRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$deadRange.Delete()

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
